# Regs_...xlsx: add two new rows describing the "discrete inputs state" and
# "current damper position" registers (30005 / 30006) to the status-register
# table, matching commit "Fix of switching to start position".

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 34: Состояние дискретных входов / 30005 / Read --------------------
$ws.Range("A34").Value = "Состояние дискретных входов"
$ws.Range("A34").WrapText = $true
$ws.Range("B34").Value = 30005
$ws.Range("C34").Value = "Read"
$ws.Rows.Item(34).RowHeight = 30

# --- Row 35: Текущее положение заслонки / 30006 / Read ---------------------
$ws.Range("A35").Value = "Текущее положение заслонки"
$ws.Range("A35").WrapText = $true
$ws.Range("B35").Value = 30006
$ws.Range("C35").Value = "Read"
$ws.Rows.Item(35).RowHeight = 30

# --- Update the view so the new rows are visible / selected ----------------
$win = $excel.ActiveWindow
$win.ScrollRow = 19
$win.ScrollColumn = 1
$ws.Range("B35").Select()

$wb.Save()
